$wb = $excel.ActiveWorkbook

# Sheets in this workbook (tab order): DATA_ELEMENT_A, DATA_ELEMENT_B, DATA_ELEMENT_C
$wsA = $wb.Worksheets.Item("DATA_ELEMENT_A")
$wsB = $wb.Worksheets.Item("DATA_ELEMENT_B")
$wsC = $wb.Worksheets.Item("DATA_ELEMENT_C")

# --- Add the new "Project Code (attribute option)" column header (G1) and its
# attribute-option UID value (H1) to every data-element sheet. G1 should pick
# up the same formatting (wrapped header style) already used by D1/E1 on each
# sheet, so copy that formatting across before writing the value.
foreach ($ws in @($wsA, $wsB, $wsC)) {
    $ws.Range("D1").Copy() | Out-Null
    $ws.Range("G1").PasteSpecial(-4122) | Out-Null
    $ws.Range("G1").Value = "Project Code (attribute option)"
    $ws.Range("H1").Value = "wr5HhbHBYfh"
}

# --- Narrow column A now that it no longer needs to fit the old long label.
$wsA.Columns.Item(1).ColumnWidth = 19.2
$wsB.Columns.Item(1).ColumnWidth = 18.5
$wsC.Columns.Item(1).ColumnWidth = 18.5

# --- Update each sheet's selection to the newly-added G1:H1 range.
$wsA.Range("G1:H1").Select() | Out-Null
$wsB.Range("G1:H1").Select() | Out-Null
$wsC.Range("G1:H1").Select() | Out-Null

# --- Make DATA_ELEMENT_C the active tab (was DATA_ELEMENT_B).
$wsC.Activate() | Out-Null
$wsC.Range("G1:H1").Select() | Out-Null
